# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
# These values were refreshed from the live data source (see commit message:
# "Update gh-pages to output generated at 456a3b4"), so only column F cells
# on the two sheets that carry the full listing need bumping to their new
# counts. The "演出" and "本地生活" sheets are untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 620
$ws1.Range("F6").Value  = 14306
$ws1.Range("F7").Value  = 16379
$ws1.Range("F9").Value  = 93
$ws1.Range("F10").Value = 4
$ws1.Range("F12").Value = 198
$ws1.Range("F21").Value = 1250
$ws1.Range("F24").Value = 35
$ws1.Range("F26").Value = 6622
$ws1.Range("F27").Value = 969
$ws1.Range("F29").Value = 17
$ws1.Range("F30").Value = 1114
$ws1.Range("F32").Value = 5718
$ws1.Range("F35").Value = 182
$ws1.Range("F36").Value = 4770

# --- Sheet "全部类型" --------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 620
$ws4.Range("F6").Value  = 14306
$ws4.Range("F7").Value  = 16379
$ws4.Range("F9").Value  = 93
$ws4.Range("F10").Value = 4
$ws4.Range("F12").Value = 198
$ws4.Range("F21").Value = 1250
$ws4.Range("F25").Value = 35
$ws4.Range("F27").Value = 6622
$ws4.Range("F28").Value = 969
$ws4.Range("F30").Value = 17
$ws4.Range("F31").Value = 1114
$ws4.Range("F35").Value = 5718
$ws4.Range("F38").Value = 182
$ws4.Range("F39").Value = 4770
